$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.980.36"
$ws.Range("E2").Value = "  -1.37%  "

$ws.Range("D3").Value = "1.818.54"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4669"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3650"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07220"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8591"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.75"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.28%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.867.42"
$ws.Range("E12").Value = "  -1.60%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07535"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.80%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.323"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.77"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.48%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.481"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.61%  "

$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008622"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.47"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.18%  "

$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "26.600.11"
$ws.Range("E21").Value = "  -3.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.146"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.51"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.48%  "

$ws.Range("D24").Value = "2.099.33"
$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.64"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.838"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.06"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.057"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.08%  "

$ws.Range("E29").Value = "  -2.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.55"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08875"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.965"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.414"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.85%  "

$ws.Range("E34").Value = "  -4.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7159"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.079"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05246"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.90%  "

$ws.Range("E38").Value = "  -1.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.919"
$ws.Range("D39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.367"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.136"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5146"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8595"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -14.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1623"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.141"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4808"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.08"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.64"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.15%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.618"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.31%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06238"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.99%  "
